$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 7-31, leaving the 5 collapsed rows for cards (rows 2-6)
$ws.Range("A7:A31").EntireRow.Delete()

$ws.Range("A2").Value = '(''Gavi, Nest Warden'', [''{2}{U}{R}{W}'', ''Legendary Creature — Human Shaman'', ''You may pay {0} rather than pay the cycling cost of the first card you cycle each turn.'', ''Whenever you draw your second card each turn, create a 2/2 red and white Dinosaur Cat creature token.'', ''2/5''])'
$ws.Range("A3").Value = '(''Jirina Kudro'', [''{1}{R}{W}{B}'', ''Legendary Creature — Human Soldier'', ''When Jirina Kudro enters the battlefield, create a 1/1 white Human Soldier creature token for each time you’ve cast a commander from the command zone this game.'', ''Other Humans you control get +2/+0.'', ''3/3''])'
$ws.Range("A4").Value = '(''Kalamax, the Stormsire'', [''{1}{G}{U}{R}'', ''Legendary Creature — Elemental Dinosaur'', ''Whenever you cast your first instant spell each turn, if Kalamax, the Stormsire is tapped, copy that spell. You may choose new targets for the copy.'', ''Whenever you copy an instant spell, put a +1/+1 counter on Kalamax.'', ''4/4''])'
$ws.Range("A5").Value = '(''Kathril, Aspect Warper'', [''{2}{W}{B}{G}'', ''Legendary Creature — Nightmare Insect'', ''When Kathril, Aspect Warper enters the battlefield, put a flying counter on any creature you control if a creature card in your graveyard has flying. Repeat this process for first strike, double strike, deathtouch, hexproof, indestructible, lifelink, menace, reach, trample, and vigilance. Then put a +1/+1 counter on Kathril for each counter put on a creature this way.'', ''3/3''])'
$ws.Range("A6").Value = '(''Otrimi, the Ever-Playful'', [''{3}{B}{G}{U}'', ''Legendary Creature — Nightmare Beast'', ''Mutate {1}{B}{G}{U} (If you cast this spell for its mutate cost, put it over or under target non-Human creature you own. They mutate into the creature on top plus all abilities from under it.)'', ''Trample'', ''Whenever this creature deals combat damage to a player, return target creature card with mutate from your graveyard to your hand.'', ''6/6''])'
